$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("en")

# Insert a new row above current row 10 (the "mineral" row that was previously there
# shifts down to row 11, etc.). This mirrors the "magma_title"/"MAGMA" row (row 9)
# pattern: a new section-title row "minerals" / "Minerals" is added right after it.
$ws.Rows.Item(10).Insert()

$ws.Cells.Item(10, 1).Value = "minerals"
$ws.Cells.Item(10, 2).Value = "Minerals"

# Update the selection to match the post-edit state (activeCell now B10).
$ws.Range("B10").Select()
